$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the authoritative diff
$updates = [ordered]@{
    'D2' = '26.332.48'
    'E2' = '  -2.21%  '
    'D3' = '1.793.01'
    'E3' = '  -2.11%  '
    'D4' = '1.004'
    'E4' = '  -0.21%  '
    'D5' = '1.004'
    'E5' = '  -0.33%  '
    'D6' = '306.70'
    'E6' = '  -1.37%  '
    'D7' = '0.4505'
    'E7' = '  -1.56%  '
    'D8' = '0.3594'
    'E8' = '  -2.71%  '
    'D9' = '45.85'
    'E9' = '  +0.05%  '
    'E10' = '  -1.56%  '
    'D11' = '0.8832'
    'E11' = '  +0.69%  '
    'D12' = '0.07727'
    'E12' = '  -1.11%  '
    'D13' = '19.37'
    'E13' = '  -1.43%  '
    'D14' = '1.821.61'
    'E14' = '  -0.26%  '
    'D15' = '5.279'
    'E15' = '  -1.12%  '
    'D16' = '6.317'
    'E16' = '  -1.31%  '
    'D17' = '84.73'
    'E17' = '  -2.92%  '
    'D18' = '1.007'
    'E18' = '  -0.22%  '
    'D19' = '0.000008517'
    'E19' = '  -2.32%  '
    'D20' = '1.004'
    'E20' = '  -0.27%  '
    'D21' = '14.25'
    'E21' = '  -1.82%  '
    'D22' = '26.357.60'
    'E22' = '  -2.22%  '
    'D23' = '4.979'
    'E23' = '  -0.59%  '
    'D24' = '10.53'
    'E24' = '  +0.89%  '
    'D25' = '2.003.44'
    'E25' = '  -2.29%  '
    'D26' = '1.971'
    'E26' = '  -2.62%  '
    'D27' = '151.10'
    'E27' = '  -0.25%  '
    'D28' = '17.82'
    'E28' = '  -2.22%  '
    'D29' = '2.020'
    'E29' = '  +2.60%  '
    'D30' = '111.93'
    'E30' = '  -1.89%  '
    'D31' = '4.881'
    'E31' = '  -1.22%  '
    'D32' = '0.08672'
    'E32' = '  -1.49%  '
    'D33' = '3.056'
    'E33' = '  +0.65%  '
    'D34' = '2.727'
    'E34' = '  +6.13%  '
    'D35' = '4.437'
    'E35' = '  -1.01%  '
    'D36' = '0.7243'
    'E36' = '  -4.04%  '
    'D37' = '1.103'
    'E37' = '  -2.81%  '
    'D38' = '1.003'
    'E38' = '  -0.09%  '
    'D39' = '1.066'
    'E39' = '  -2.14%  '
    'E40' = '  -0.23%  '
    'D41' = '0.05083'
    'E41' = '  -1.49%  '
    'D42' = '2.854'
    'E42' = '  -1.08%  '
    'D43' = '6.857'
    'E43' = '  -1.32%  '
    'D44' = '0.5053'
    'E44' = '  +1.39%  '
    'E45' = '  -5.66%  '
    'D46' = '7.991'
    'E46' = '  -3.85%  '
    'D47' = '1.004'
    'E47' = '  -0.41%  '
    'D48' = '0.4625'
    'E48' = '  -1.38%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D49' = '9.887'
    'E49' = '  -2.31%  '
    'B50' = 'Quant'
    'C50' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D50' = '100.87'
    'E50' = '  -1.47%  '
    'D51' = '1.579'
    'E51' = '  -2.24%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($addr[0] -eq "D") {
        # Price column: values look numeric (e.g. "1.004", "306.70") but must
        # stay plain text exactly as scraped, so force text format before writing
        # and then drop back to the default style so no stray formatting remains.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
